$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text that often
# *looks* numeric ("0.999", "63.910.61", "  -1.00%  "). Flip the whole
# data range to text format before writing so COM stores these as
# strings (matching the workbook's original inlineStr cells) instead of
# coercing them into real numbers, then restore the Normal style so no
# stray number format is left on the cells.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.910.61"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "3.063.37"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "559.95"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "142.73"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "3.061.57"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "35.42"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "3.561.69"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "63.924.53"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "3.062.37"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "488.20"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").Value = "0.692"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  +8.17%  "
$ws.Range("D24").Value = "7.53"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "82.68"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "8.21"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "26.49"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "2.59"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("D35").Value = "6.27"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").Value = "54.98"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "0.0412"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "445.19"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "0.0816"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("D40").Value = "3.045.79"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  -8.37%  "
$ws.Range("D42").Value = "8.35"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("D45").Value = "28.29"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("D48").Value = "0.114"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "117.66"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  +3.78%  "

$priceVolRange.Style = "Normal"
